# Apply "Add files via upload: added Ixz, compressive stress, tensile stress"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Range("F1").Value = "CompressiveStress"
$ws.Range("G1").Value = "TensileStress"
$ws.Range("H1").Value = "Ix"
$ws.Range("I1").Value = "Iy"
$ws.Range("J1").Value = "Iz"
$ws.Range("K1").Value = "Ixy"
$ws.Range("L1").Value = "Iyz"
$ws.Range("M1").Value = "Ixz"

# ---- Shared "Tip Chord" formula across C2:C7 (keeps the existing C2:C4 shared group intact) ----
$ws.Range("C2:C7").Formula = "=8.4"

# ---- Row 2 ----
$ws.Range("A2").Value = 40
$ws.Range("B2").Value = 12.6
$ws.Range("D2").Formula = "=1811*2"
$ws.Range("E2").Value = 2.097
$ws.Range("F2").Formula = "=1.818*10^7"
$ws.Range("G2").Formula = "=3.008*10^8"
$ws.Range("H2").Formula = "=2420160"
$ws.Range("I2").Formula = "=2589780"
$ws.Range("J2").Formula = "=177178"
$ws.Range("K2").Formula = "=34.3927"
$ws.Range("L2").Formula = "=-639.485"
$ws.Range("M2").Formula = "=129704"

# ---- Row 3 ----
$ws.Range("A3").Value = 37
$ws.Range("B3").Value = 12.6
$ws.Range("D3").Formula = "=1676*2"
$ws.Range("E3").Value = 1.581
$ws.Range("F3").Formula = "=1.54*10^7"
$ws.Range("G3").Formula = "=2.689*10^8"
$ws.Range("H3").Value = 1919140
$ws.Range("I3").Formula = "=2076080"
$ws.Range("J3").Formula = "=163938"
$ws.Range("K3").Formula = "=31.8494"
$ws.Range("L3").Formula = "=-548.076"
$ws.Range("M3").Formula = "=111080"

# ---- Row 4 ----
$ws.Range("A4").Value = 35
$ws.Range("B4").Value = 12.6
$ws.Range("D4").Formula = "=1586*2"
$ws.Range("E4").Value = 1.296
$ws.Range("F4").Formula = "=1.399*10^7"
$ws.Range("G4").Formula = "=2.484*10^8"
$ws.Range("H4").Formula = "=1626860"
$ws.Range("I4").Formula = "=1775350"
$ws.Range("J4").Formula = "=155112"
$ws.Range("K4").Formula = "=30.1543"
$ws.Range("L4").Formula = "=-491.052"
$ws.Range("M4").Formula = "=99466.6"

# ---- Row 5 (new row) ----
$ws.Range("A5").Value = 33
$ws.Range("B5").Value = 12.6
$ws.Range("D5").Formula = "=1496*2"
$ws.Range("E5").Value = 1.052
$ws.Range("F5").Formula = "=1.275*10^7"
$ws.Range("G5").Formula = "=2.285*10^8"
$ws.Range("H5").Formula = "=1365890"
$ws.Range("I5").Formula = "=1505940"
$ws.Range("J5").Formula = "=146286"
$ws.Range("K5").Formula = "=28.4596"
$ws.Range("L5").Formula = "=-437.16"
$ws.Range("M5").Formula = "=88494.2"

# ---- Row 6 ----
$ws.Range("A6").Value = 30
$ws.Range("B6").Value = 12.6
$ws.Range("D6").Formula = "=1361*2"
$ws.Range("E6").Value = 0.7551
$ws.Range("F6").Formula = "=1.072*10^7"
$ws.Range("G6").Formula = "=1.996*10^8"
$ws.Range("H6").Formula = "=1029270"
$ws.Range("I6").Formula = "=1156640"
$ws.Range("J6").Formula = "=133049"
$ws.Range("K6").Formula = "=25.9188"
$ws.Range("L6").Formula = "=-362.196"
$ws.Range("M6").Formula = "=73238.6"

# ---- Row 7 ----
$ws.Range("A7").Value = 20
$ws.Range("B7").Value = 12.6
$ws.Range("D7").Formula = "=911*2"
$ws.Range("E7").Value = 0.02218
$ws.Range("F7").Formula = "=5.942*10^6"
$ws.Range("G7").Formula = "=1.351*10^8"
$ws.Range("H7").Formula = "=310232"
$ws.Range("I7").Formula = "=395391"
$ws.Range("J7").Formula = "=88948.3"
$ws.Range("K7").Formula = "=17.4658"
$ws.Range("L7").Formula = "=-163.216"
$ws.Range("M7").Formula = "=32811.7"

# ---- Column widths ----
$ws.Columns.Item(5).ColumnWidth = 11.7109375
$ws.Columns.Item(6).ColumnWidth = 18.7109375
$ws.Columns.Item(7).ColumnWidth = 18.7109375
$ws.Columns.Item(8).ColumnWidth = 14.140625
$ws.Columns.Item(9).ColumnWidth = 9.28515625
$ws.Columns.Item(10).ColumnWidth = 9.28515625
$ws.Columns.Item(11).ColumnWidth = 9.28515625
$ws.Columns.Item(12).ColumnWidth = 9.28515625

# ---- Number format: F loses the old scientific style, H2:H3 gain it ----
$ws.Range("F2:F7").ClearFormats()
$ws.Range("H2:H3").NumberFormat = "0.00E+00"

# ---- Selection ----
$ws.Range("O7").Select()
